$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.725.34'
$ws.Range("E2").Value = '  +0.00%  '
$ws.Range("D3").Value = '3.311.14'
$ws.Range("E3").Value = '  -1.87%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '579.89'
$ws.Range("E5").Value = '  -2.19%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '173.22'
$ws.Range("E6").Value = '  -7.37%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("E8").Value = '  -2.81%  '
$ws.Range("D9").Value = '3.306.88'
$ws.Range("E9").Value = '  -1.86%  '
$ws.Range("E10").Value = '  -4.93%  '
$ws.Range("E11").Value = '  -2.72%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '45.24'
$ws.Range("E12").Value = '  -4.66%  '
$ws.Range("E13").Value = '  -3.35%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '661.36'
$ws.Range("E14").Value = '  +3.67%  '
$ws.Range("D15").Value = '3.848.97'
$ws.Range("E15").Value = '  -1.73%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '8.35'
$ws.Range("E16").Value = '  -3.24%  '
$ws.Range("D17").Value = '67.848.83'
$ws.Range("E17").Value = '  +0.14%  '
$ws.Range("E18").Value = '  -1.03%  '
$ws.Range("D19").Value = '3.314.76'
$ws.Range("E19").Value = '  -2.00%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.38'
$ws.Range("E20").Value = '  -3.91%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.85'
$ws.Range("E21").Value = '  -2.79%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.885'
$ws.Range("E22").Value = '  -2.78%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.42'
$ws.Range("E23").Value = '  +6.03%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '16.86'
$ws.Range("E24").Value = '  -5.99%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '96.89'
$ws.Range("E25").Value = '  -3.09%  '
$ws.Range("E26").Value = '  -5.30%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.66'
$ws.Range("E27").Value = '  -7.07%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.26'
$ws.Range("E28").Value = '  -5.16%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '33.27'
$ws.Range("E29").Value = '  +1.23%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.36'
$ws.Range("E30").Value = '  -4.13%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.25'
$ws.Range("E31").Value = '  +5.00%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '589.86'
$ws.Range("E32").Value = '  -3.82%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '10.91'
$ws.Range("E33").Value = '  -1.78%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.103'
$ws.Range("E34").Value = '  -2.73%  '
$ws.Range("B35").Value = 'Dai'
$ws.Range("C35").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  -0.10%  '
$ws.Range("B36").Value = 'Maker'
$ws.Range("C36").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D36").Value = '3.708.61'
$ws.Range("E36").Value = '  -8.47%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '56.86'
$ws.Range("E37").Value = '  +0.90%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.27'
$ws.Range("E38").Value = '  -14.44%  '
$ws.Range("E39").Value = '  -0.83%  '
$ws.Range("B40").Value = 'InjectiveProtocol'
$ws.Range("C40").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '32.34'
$ws.Range("E40").Value = '  -4.79%  '
$ws.Range("B41").Value = 'Fetch.AI'
$ws.Range("C41").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.60'
$ws.Range("E41").Value = '  -6.92%  '
$ws.Range("E42").Value = '  -5.17%  '
$ws.Range("B43").Value = 'ApeXProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.29'
$ws.Range("E43").Value = '  -3.77%  '
$ws.Range("B44").Value = 'TheGraph'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.330'
$ws.Range("E44").Value = '  -3.82%  '
$ws.Range("D45").Value = '0.0₃0659'
$ws.Range("E45").Value = '  -5.96%  '
$ws.Range("E46").Value = '  -3.78%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.58'
$ws.Range("E47").Value = '  -0.55%  '
$ws.Range("E48").Value = '  -2.23%  '
$ws.Range("E50").Value = '  -3.95%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '127.50'
$ws.Range("E51").Value = '  -0.08%  '
